# Insert two new weekly records for "Pimiento" (Zafiro rojo / Zafiro verde)
# at the top of the Terminal Hortofrutícola Agro Chillán price history
# block, pushing the existing rows 129-147 down to 131-149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("129:130").Insert()

# Row 129 - Zafiro rojo
$ws.Range("A129").Value = 7
$ws.Range("B129").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C129").Value = "Ñuble"
$ws.Range("D129").Value = 44474
$ws.Range("E129").Value = 16
$ws.Range("F129").Value = 100112002
$ws.Range("G129").Value = "Pimiento"
$ws.Range("H129").Value = "Zafiro rojo"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 60
$ws.Range("K129").Value = 42000
$ws.Range("L129").Value = 43000
$ws.Range("M129").Value = 42500
$ws.Range("N129").Value = "`$/caja 15 kilos"
$ws.Range("O129").Value = "Región de Arica y Parinacota"
$ws.Range("P129").Value = 2833
$ws.Range("Q129").Value = 15
$ws.Range("R129").Value = "Hortaliza"

# Row 130 - Zafiro verde
$ws.Range("A130").Value = 7
$ws.Range("B130").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C130").Value = "Ñuble"
$ws.Range("D130").Value = 44474
$ws.Range("E130").Value = 16
$ws.Range("F130").Value = 100112002
$ws.Range("G130").Value = "Pimiento"
$ws.Range("H130").Value = "Zafiro verde"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 60
$ws.Range("K130").Value = 38000
$ws.Range("L130").Value = 39000
$ws.Range("M130").Value = 38500
$ws.Range("N130").Value = "`$/caja 15 kilos"
$ws.Range("O130").Value = "Región de Arica y Parinacota"
$ws.Range("P130").Value = 2567
$ws.Range("Q130").Value = 15
$ws.Range("R130").Value = "Hortaliza"
